# Ajuste validacion amortizacion X Naturaleza de cuentas
# Update the test-data row on the "AplicacionPago" sheet with the new
# pagaduria / period values used by the automation test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data row (row 2) -------------------------------------------------

# IdPagaduria: 271 -> 103 (now stored as text, matching the column's
# text number format)
$ws.Range("A2").Value = "103"

# Periodo: Octubre 30 -> Noviembre 30
$ws.Range("B2").Value = "Noviembre 30"

# NombrePagaduria
$ws.Range("C2").Value = '"CONSORCIO DE PENSIONES DEL HUILA"'

# RutaPagaduria (unchanged, left as-is)

# Ano (unchanged, left as-is)

# PeriodoEspacio: "Octubre  30" -> "Noviembre  30"
$ws.Range("F2").Value = '"Noviembre  30"'

# FiltroFecha: 30/10/2021 -> 30/11/2021
$ws.Range("G2").Value = "30/11/2021"

# AccountingSource (unchanged, left as-is)

# AccountingName: added the second "upper(...)" clause
$ws.Range("I2").Value = "`"upper('Aplicación de pago por pagaduría'),  upper('Aplicación de pago venta en firme') `""

# FechaRegistro: 25/11/2021 -> 10/12/2021
$ws.Range("J2").Value = "10/12/2021"

# --- number formats / style clean-up -----------------------------------
# These cells keep their text ("@") format; re-applying it lets identical
# xf records collapse instead of keeping the old, now-unused duplicates.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"

# H2 ("'APLPAG'") previously carried a quote-prefix style; clear that
# formatting quirk and restore the plain text format + value.
$ws.Range("H2").ClearFormats()
$ws.Range("H2").Value = "`"'APLPAG'`""
$ws.Range("H2").NumberFormat = "@"

# Match the saved cursor/selection position from the authoring session.
[void]$ws.Range("C17").Select()
